$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.369.92'
$ws.Range('E2').Value = '  +4.42%  '
$ws.Range('D3').Value = '2.046.96'
$ws.Range('E3').Value = '  +2.88%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.654'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '65.96'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +10.89%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.410'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +12.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0787'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.09%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.930'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +26.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('D16').Value = '2.346.87'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.70%  '
$ws.Range('D18').Value = '2.048.99'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').Value = '37.276.24'
$ws.Range('E19').Value = '  +4.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.56%  '
$ws.Range('D21').Value = '0.0₃0888'
$ws.Range('E21').Value = '  +4.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.56'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.64'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.52%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +4.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.135'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +36.69%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.87%  '
$ws.Range('E31').Value = '  +2.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.59%  '
$ws.Range('E33').Value = '  +5.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0634'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.00%  '
$ws.Range('E35').Value = '  +6.97%  '
$ws.Range('E36').Value = '  -3.36%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.59%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  +3.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +31.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.102'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.79%  '
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('E43').Value = '  +6.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '17.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.97%  '
$ws.Range('E46').Value = '  +2.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '96.73'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('D49').Value = '1.407.42'
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.63%  '
